# Implementação da Função calcular_beneficios_inss
#
# Adds the new "calcular_beneficios_inss" function to the Funcoes_Inputs and
# Funcoes_Outputs sheets (mirroring the existing calcular_* blocks), and
# updates the active-sheet / selection bookkeeping to match where the author
# left the cursor after the edit.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Funcoes_Inputs: new rows 45-60 for calcular_beneficios_inss inputs
# ---------------------------------------------------------------------
$wsIn = $wb.Worksheets.Item("Funcoes_Inputs")

$inputsData = @(
    @("Nev_Safast_Tipico"),
    @("Nev_Safast_Trajeto"),
    @("Nev_Safast_DoenOcup"),
    @("Nev_Safast_NRelac"),
    @("Nev_Obito_Tipico"),
    @("Nev_Obito_Trajeto"),
    @("Nev_Obito_DoenOcup"),
    @("Nev_Obito_NRelac"),
    @("Nev_Afmenor15_Tipico"),
    @("Nev_Afmenor15_Trajeto"),
    @("Nev_Afmenor15_DoenOcup"),
    @("Nev_Afmenor15_NRelac"),
    @("Nev_Afmaior15_Tipico"),
    @("Nev_Afmaior15_Trajeto"),
    @("Nev_Afmaior15_DoenOcup"),
    @("Nev_Afmaior15_NRelac")
)

$row = 45
foreach ($item in $inputsData) {
    $wsIn.Cells.Item($row, 1).Value = "calcular_beneficios_inss"
    $wsIn.Cells.Item($row, 2).Value = $item[0]
    $row = $row + 1
}

# ---------------------------------------------------------------------
# Funcoes_Outputs: new rows 29-34 for calcular_beneficios_inss outputs
# ---------------------------------------------------------------------
$wsOut = $wb.Worksheets.Item("Funcoes_Outputs")

$outputsData = @(
    @("NB_91"),
    @("NB_92"),
    @("NB_93"),
    @("NB_94"),
    @("NB_31"),
    @("NB_32")
)

$row = 29
foreach ($item in $outputsData) {
    $wsOut.Cells.Item($row, 1).Value = "calcular_beneficios_inss"
    $wsOut.Cells.Item($row, 2).Value = $item[0]
    $row = $row + 1
}

# ---------------------------------------------------------------------
# Cursor / selection bookkeeping left behind by the edit
# ---------------------------------------------------------------------

# Parametros was the previously-active sheet; its selection moved on (but it
# is no longer the active tab).
$wsParam = $wb.Worksheets.Item("Parametros")
$wsParam.Range("D35").Select()

# Funcoes_Inputs selection after appending the new rows.
$wsIn.Range("A45").Select()

# Funcoes_Outputs ends up as the active sheet/tab with the cursor parked
# after the newly-added rows.
$wsOut.Activate()
$wsOut.Range("B36").Select()
